$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate labels first (B14 "ID", B15 "PASS"), mirroring rows 11/12 ---
$ws.Range("B14").Value = "ID"
$ws.Range("B11").Copy()
$ws.Range("B14").PasteSpecial(-4122)

$ws.Range("B15").Value = "PASS"
$ws.Range("B12").Copy()
$ws.Range("B15").PasteSpecial(-4122)

# --- Row 15 (PASS / password123) styled first so its new font/xf is registered
#     before the hyperlink-wrap xf, matching the order the workbook was edited in ---
$ws.Range("C15").Value = "password123"
$pf = $ws.Range("C15").Font
$pf.Name = "Courier New"
$pf.Family = 3
$pf.Size = 11
$pf.Color = 1381795
$ws.Range("C15").WrapText = $true
$ws.Range("C15").HorizontalAlignment = -4131
$ws.Range("C15").VerticalAlignment = -4108
$ws.Range("C15").IndentLevel = 1

# --- Row 14 (ID / email hyperlink), mirrors row 11's ID/URL row ---
$ws.Range("C14").Value = "shinji19750918@yahoo.co.jp"
$ws.Hyperlinks.Add($ws.Range("C14"), "mailto:shinji19750918@yahoo.co.jp")
$ws.Range("C14").WrapText = $true
$ws.Range("C14").HorizontalAlignment = -4131
$ws.Range("C14").VerticalAlignment = -4108
$ws.Range("C14").IndentLevel = 1

# --- Row heights (rows 2-5 got a slightly taller default in the resave) ---
$ws.Rows.Item(2).RowHeight = 20.25
$ws.Rows.Item(3).RowHeight = 20.25
$ws.Rows.Item(4).RowHeight = 20.25
$ws.Rows.Item(5).RowHeight = 20.25
$ws.Rows.Item(14).RowHeight = 59.65
$ws.Rows.Item(15).RowHeight = 28.5

# --- Column C width to fit the new content ---
$ws.Columns.Item(3).ColumnWidth = 57.48

# --- Selection moves to the newly added password cell ---
$ws.Range("C15").Select()

Write-Output "edit complete"
